$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1246.75
$ws.Range("I40").Value = 814.75
$ws.Range("J40").Value = 1462.75
$ws.Range("K40").Value = 814.75
$ws.Range("L40").Value = 1462.75
$ws.Range("M40").Value = -639.75
$ws.Range("N40").Value = -1812.75
$ws.Range("H129").Value = 884.4915
$ws.Range("I129").Value = 626
$ws.Range("J129").Value = 898.3393
$ws.Range("K129").Value = 1878
$ws.Range("L129").Value = 2695.0179
$ws.Range("M129").Value = 3122
$ws.Range("N129").Value = -12695.0179
$ws.Range("H137").Value = 1679.6923
$ws.Range("I137").Value = 1235.75
$ws.Range("K137").Value = 3707.25
$ws.Range("M137").Value = -1157.25
$ws.Range("H138").Value = 3073.875
$ws.Range("J138").Value = 3333.5862
$ws.Range("L138").Value = 10000.7586
$ws.Range("N138").Value = -20280.7586

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2588.9216
$ws.Range("I32").Value = 1764.262
$ws.Range("K32").Value = 1764.262
$ws.Range("M32").Value = -1477.262
$ws.Range("H45").Value = 1837.5238
$ws.Range("I45").Value = 1409.3334
$ws.Range("K45").Value = 1409.3334
$ws.Range("M45").Value = -1032.3334
$ws.Range("H74").Value = 1094.5555
$ws.Range("I74").Value = 542.4
$ws.Range("J74").Value = 1306.9231
$ws.Range("K74").Value = 542.4
$ws.Range("L74").Value = 1306.9231
$ws.Range("M74").Value = 331.6
$ws.Range("N74").Value = -3054.9231
$ws.Range("H77").Value = 1094.5555
$ws.Range("I77").Value = 542.4
$ws.Range("J77").Value = 1306.9231
$ws.Range("K77").Value = 2712
$ws.Range("L77").Value = 6534.6155
$ws.Range("M77").Value = 1656
$ws.Range("N77").Value = -15270.6155

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1483.6
$ws.Range("J86").Value = 1991.2
$ws.Range("L86").Value = 1991.2
$ws.Range("N86").Value = -4237.2
$ws.Range("H89").Value = 1483.6
$ws.Range("J89").Value = 1991.2
$ws.Range("L89").Value = 9956
$ws.Range("N89").Value = -21188

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9851.485000000001
$ws.Range("I31").Value = 12269.308
$ws.Range("J31").Value = 2866.6667
$ws.Range("K31").Value = 12269.308
$ws.Range("L31").Value = 2866.6667
$ws.Range("M31").Value = -11974.308
$ws.Range("N31").Value = -3456.6667
$ws.Range("H34").Value = 9851.485000000001
$ws.Range("I34").Value = 12269.308
$ws.Range("J34").Value = 2866.6667
$ws.Range("K34").Value = 12269.308
$ws.Range("L34").Value = 2866.6667
$ws.Range("M34").Value = -12067.308
$ws.Range("N34").Value = -3270.6667
$ws.Range("H141").Value = 44700
$ws.Range("J141").Value = 44700
$ws.Range("L141").Value = 44700
$ws.Range("N141").Value = -55060

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 782.3421
$ws.Range("I5").Value = 897.5714
$ws.Range("J5").Value = 715.125
$ws.Range("K5").Value = 2692.7142
$ws.Range("L5").Value = 2145.375
$ws.Range("M5").Value = -2580.7142
$ws.Range("N5").Value = -2369.375
$ws.Range("H92").Value = 612.5
$ws.Range("J92").Value = 557.1429000000001
$ws.Range("L92").Value = 1671.4287
$ws.Range("N92").Value = -4167.4287
$ws.Range("H97").Value = 878.9167
$ws.Range("J97").Value = 1136.5
$ws.Range("L97").Value = 3409.5
$ws.Range("N97").Value = -4401.5
$ws.Range("H131").Value = 801.1111
$ws.Range("J131").Value = 821.1579
$ws.Range("L131").Value = 2463.4737
$ws.Range("N131").Value = -12543.4737
$ws.Range("H135").Value = 782.3421
$ws.Range("I135").Value = 897.5714
$ws.Range("J135").Value = 715.125
$ws.Range("K135").Value = 8078.1426
$ws.Range("L135").Value = 6436.125
$ws.Range("M135").Value = -5543.1426
$ws.Range("N135").Value = -11506.125

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3552.111
$ws.Range("I122").Value = 2496.125
$ws.Range("J122").Value = 12000
$ws.Range("K122").Value = 7488.375
$ws.Range("L122").Value = 36000
$ws.Range("M122").Value = -5038.375
$ws.Range("N122").Value = -40900
$ws.Range("H126").Value = 5013.793
$ws.Range("I126").Value = 4077.7778
$ws.Range("J126").Value = 6545.4546
$ws.Range("K126").Value = 12233.3334
$ws.Range("L126").Value = 19636.3638
$ws.Range("M126").Value = -9763.3334
$ws.Range("N126").Value = -24576.3638

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6324.875
$ws.Range("I7").Value = 6128.4287
$ws.Range("J7").Value = 7700
$ws.Range("K7").Value = 6128.4287
$ws.Range("L7").Value = 7700
$ws.Range("M7").Value = -6016.4287
$ws.Range("N7").Value = -7924
$ws.Range("H16").Value = 935
$ws.Range("I16").Value = 935
$ws.Range("K16").Value = 935
$ws.Range("M16").Value = -765
$ws.Range("H40").Value = 3222.423
$ws.Range("I40").Value = 2980.1428
$ws.Range("J40").Value = 4240
$ws.Range("K40").Value = 2980.1428
$ws.Range("L40").Value = 4240
$ws.Range("M40").Value = -2844.1428
$ws.Range("N40").Value = -4512
$ws.Range("H46").Value = 675.43475
$ws.Range("I46").Value = 620.58826
$ws.Range("J46").Value = 830.8333
$ws.Range("K46").Value = 620.58826
$ws.Range("L46").Value = 830.8333
$ws.Range("M46").Value = -432.58826
$ws.Range("N46").Value = -1206.8333
$ws.Range("H82").Value = 1950.6923
$ws.Range("I82").Value = 3063.3333
$ws.Range("J82").Value = 997
$ws.Range("K82").Value = 3063.3333
$ws.Range("L82").Value = 997
$ws.Range("M82").Value = -2702.3333
$ws.Range("N82").Value = -1719
$ws.Range("H85").Value = 1950.6923
$ws.Range("I85").Value = 3063.3333
$ws.Range("J85").Value = 997
$ws.Range("K85").Value = 3063.3333
$ws.Range("L85").Value = 997
$ws.Range("M85").Value = -1815.3333
$ws.Range("N85").Value = -3493
$ws.Range("H122").Value = 2431.1
$ws.Range("I122").Value = 1944.4286
$ws.Range("J122").Value = 3566.6667
$ws.Range("K122").Value = 5833.2858
$ws.Range("L122").Value = 10700.0001
$ws.Range("M122").Value = -3383.2858
$ws.Range("N122").Value = -15600.0001
$ws.Range("H126").Value = 6324.875
$ws.Range("I126").Value = 6128.4287
$ws.Range("J126").Value = 7700
$ws.Range("K126").Value = 18385.2861
$ws.Range("L126").Value = 23100
$ws.Range("M126").Value = -15915.2861
$ws.Range("N126").Value = -28040
